{"js": "// Remove the trailing \"Ver no Jupiter ...\" line and the\n// \"\u00a9 2020 ... Creative Commons Attribution\" line (the Jekyll site footer),\n// plus the now-redundant blank paragraph left behind once they are gone \u2014\n// the footer text paragraphs were sandwiched between two blank paragraphs,\n// and after removing the footer text only a single blank paragraph remains\n// before the trailing page-break paragraph.\n\nconst body = context.document.body;\nlet paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet items = paragraphs.items;\n\n// Locate the two footer paragraphs by their exact text.\nlet verNoJupiterIdx = -1;\nlet copyrightIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text;\n  if (verNoJupiterIdx === -1 && text === \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n    verNoJupiterIdx = i;\n  } else if (copyrightIdx === -1 && text.indexOf(\"\u00a9 2020\") === 0) {\n    copyrightIdx = i;\n  }\n}\n\nif (verNoJupiterIdx === -1 || copyrightIdx === -1) {\n  throw new Error(\"Could not locate the footer paragraphs to remove.\");\n}\n\n// Delete the later paragraph first so the earlier index remains valid.\nconst firstIdx = Math.min(verNoJupiterIdx, copyrightIdx);\nconst secondIdx = Math.max(verNoJupiterIdx, copyrightIdx);\n\nitems[secondIdx].delete();\nitems[firstIdx].delete();\nawait context.sync();\n\n// Re-load the paragraphs collection so indices/text reflect the deletions\n// that just happened (the old `items` array still points at the now-removed\n// paragraphs).\nparagraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\nitems = paragraphs.items;\n\n// Collapse the duplicate blank paragraph: if the paragraph now sitting\n// right before where the footer used to be, and the one right after it, are\n// both empty, drop the earlier one so a single blank paragraph separates\n// the requirements list from the trailing page-break paragraph (matching\n// the source change).\nconst before = items[firstIdx - 1];\nconst after = items[firstIdx];\nif (before && after && before.text === \"\" && after.text === \"\") {\n  before.delete();\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" line and the\n# \"\u00a9 2020 ... Creative Commons Attribution\" line (the Jekyll site footer),\n# plus the now-redundant blank paragraph left behind once they are gone \u2014\n# the footer text paragraphs were sandwiched between two blank paragraphs,\n# and after removing the footer text only a single blank paragraph remains\n# before the trailing page-break paragraph.\n\n$d = $word.ActiveDocument\n\n# Locate the two footer paragraphs by their exact text.\n$verIdx = -1\n$copyIdx = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)\n    if ($verIdx -eq -1 -and $t -eq \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n        $verIdx = $i\n    } elseif ($copyIdx -eq -1 -and $t.StartsWith(\"\u00a9 2020\")) {\n        $copyIdx = $i\n    }\n}\n\nif ($verIdx -eq -1 -or $copyIdx -eq -1) {\n    throw \"Could not locate the footer paragraphs to remove.\"\n}\n\n# Delete the later paragraph first so the earlier index stays valid.\n$firstIdx = [Math]::Min($verIdx, $copyIdx)\n$secondIdx = [Math]::Max($verIdx, $copyIdx)\n\n$d.Paragraphs.Item($secondIdx).Range.Delete()\n$d.Paragraphs.Item($firstIdx).Range.Delete()\n\n# Collapse the now-duplicate blank paragraph: the footer text used to sit\n# between two blank paragraphs, so once it is gone only a single blank\n# paragraph should remain before the trailing page-break paragraph.\n$beforeText = $d.Paragraphs.Item($firstIdx - 1).Range.Text.TrimEnd([char]13)\n$afterText = $d.Paragraphs.Item($firstIdx).Range.Text.TrimEnd([char]13)\n\nif ($beforeText -eq \"\" -and $afterText -eq \"\") {\n    $d.Paragraphs.Item($firstIdx - 1).Range.Delete()\n}\n"}
